$d = $word.ActiveDocument

function Get-ParaIndexOfText($searchText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.Paragraphs.Item(1).Index
}

# ------------------------------------------------------------------
# 1) Remove one of the extra empty (bold) paragraphs that sits right
#    before the centered "CONTACT" paragraph.
# ------------------------------------------------------------------
$contactIdx = Get-ParaIndexOfText("CONTACT")
$d.Paragraphs.Item($contactIdx - 1).Range.Delete()

# ------------------------------------------------------------------
# 2) Update the street-address paragraph:
#      "60 South 600 East Ste. 250"  ->  "15 West South Temple Ste. 1630"
#    split across 4 runs exactly like the authored edit:
#      "15 West South Temple" / " Ste. " / "163" / "0"
# ------------------------------------------------------------------
$addrIdx = Get-ParaIndexOfText("60 South 600 East Ste. 250")
$addrPara = $d.Paragraphs.Item($addrIdx)
$addrPara.Range.Find.Execute("60 South 600 East Ste. 250", $false, $false, $false, $false, $false, $true, 1, $false, "15 West South Temple Ste. 1630", 2) | Out-Null

$addrRange = $d.Paragraphs.Item($addrIdx).Range
$addrStart = $addrRange.Start
$addrLen = 30   # length of "15 West South Temple Ste. 1630"

foreach ($boundary in 20, 26, 29) {
    $suffix = $d.Range($addrStart + $boundary, $addrStart + $addrLen)
    $suffix.Font.Bold = $true
    $suffix.Font.Bold = $false
}

# ------------------------------------------------------------------
# 3) Update the city/state/zip paragraph and append two brand-new
#    paragraphs after it: the phone number, then a blank paragraph.
#      "Salt Lake City, Utah 84102" -> "Salt Lake City, Utah 84101"
#    (split into "Salt Lake City, Utah 8410" / "1")
#    New paragraph: "801-406-7877" (KeepTogether / w:keepLines)
#    New paragraph: empty
# ------------------------------------------------------------------
$cityIdx = Get-ParaIndexOfText("Salt Lake City, Utah 84102")
$cityPara = $d.Paragraphs.Item($cityIdx)
$cityPara.Range.Find.Execute("Salt Lake City, Utah 84102", $false, $false, $false, $false, $false, $true, 1, $false, "Salt Lake City, Utah 84101^p801-406-7877^p", 2) | Out-Null

$cityRange = $d.Paragraphs.Item($cityIdx).Range
$cityStart = $cityRange.Start
$cityLen = 26  # length of "Salt Lake City, Utah 84101"
$citySuffix = $d.Range($cityStart + 25, $cityStart + $cityLen)
$citySuffix.Font.Bold = $true
$citySuffix.Font.Bold = $false

# Phone-number paragraph (the one right after the city paragraph): split into
# "801-" / "406-7877" and mark KeepTogether (w:keepLines)
$phonePara = $d.Paragraphs.Item($cityIdx + 1)
$phonePara.KeepTogether = $true
$phoneRange = $phonePara.Range
$phoneStart = $phoneRange.Start
$phoneLen = 12  # length of "801-406-7877"
$phoneSuffix = $d.Range($phoneStart + 4, $phoneStart + $phoneLen)
$phoneSuffix.Font.Bold = $true
$phoneSuffix.Font.Bold = $false

Write-Output ("Done. Paragraphs.Count=" + $d.Paragraphs.Count)
